# Allocate revenues from transmission construction by ISIC Code (#120)
#
# Inserts a new row into the "Key to Variables" sheet (row 92) for the new
# acronym "SoTCCbIC" ("Share of Transmission Capital Costs by ISIC Code"),
# pushing all subsequent rows (old 92..221) down by one (new 93..222).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new blank row above the existing row 92 ("SYC" / Start Year
# Capacities), shifting it (and everything below it) down by one.
$ws.Rows.Item(92).Insert()

# Populate the new row with the new variable entry.
$ws.Cells.Item(92, 1).Value = "elec"
$ws.Cells.Item(92, 2).Value = "SoTCCbIC"
$ws.Cells.Item(92, 3).Value = "Share of Transmission Capital Costs by ISIC Code"
$ws.Cells.Item(92, 4).Clear()
$ws.Cells.Item(92, 6).Value = "low"
